$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value2 = 66.47695399999999
$ws.Cells.Item(2, 8).Value2 = 199.430862
$ws.Cells.Item(2, 9).Value2 = 0.04311983106164722
$ws.Cells.Item(2, 10).Value2 = 0.04311983106164721
$ws.Cells.Item(2, 11).Value2 = 3.0
$ws.Cells.Item(2, 12).Value2 = 1.0
$ws.Cells.Item(2, 13).Value2 = 0.5531836666666666
$ws.Cells.Item(2, 14).Value2 = 1.659551
$ws.Cells.Item(2, 15).Value2 = 0.05946633586794156
$ws.Cells.Item(2, 16).Value2 = 0.05946633586794157
$ws.Cells.Item(2, 17).Value2 = 36.77396516255133
$ws.Cells.Item(2, 18).Value2 = 330.965686462962
$ws.Cells.Item(2, 19).Value2 = 0.002564178356480813
$ws.Cells.Item(2, 20).Value2 = 0.002564178356480813

$ws.Cells.Item(3, 7).Value2 = 66.47695399999999
$ws.Cells.Item(3, 8).Value2 = 199.430862
$ws.Cells.Item(3, 9).Value2 = 0.04311983106164722
$ws.Cells.Item(3, 10).Value2 = 0.04311983106164721
$ws.Cells.Item(3, 11).Value2 = 3
$ws.Cells.Item(3, 12).Value2 = 1
$ws.Cells.Item(3, 13).Value2 = 4.444398333333333
$ws.Cells.Item(3, 14).Value2 = 13.333195
$ws.Cells.Item(3, 15).Value2 = 0.4777655233631019
$ws.Cells.Item(3, 16).Value2 = 0.4777655233631019
$ws.Cells.Item(3, 17).Value2 = 295.4500635626766
$ws.Cells.Item(3, 18).Value2 = 2659.05057206409
$ws.Cells.Item(3, 19).Value2 = 0.02060116865449642
$ws.Cells.Item(3, 20).Value2 = 0.02060116865449642

$ws.Cells.Item(4, 7).Value2 = 66.47695399999999
$ws.Cells.Item(4, 8).Value2 = 199.430862
$ws.Cells.Item(4, 9).Value2 = 0.04311983106164722
$ws.Cells.Item(4, 10).Value2 = 0.04311983106164721
$ws.Cells.Item(4, 11).Value2 = 2
$ws.Cells.Item(4, 12).Value2 = 0.6666666666666666
$ws.Cells.Item(4, 13).Value2 = 0.08877066666666666
$ws.Cells.Item(4, 14).Value2 = 0.266312
$ws.Cells.Item(4, 15).Value2 = 0.009542700909862518
$ws.Cells.Item(4, 16).Value2 = 0.00954270090986252
$ws.Cells.Item(4, 17).Value2 = 5.901203524549333
$ws.Cells.Item(4, 18).Value2 = 53.110831720944
$ws.Cells.Item(4, 19).Value2 = 0.000411479651105099
$ws.Cells.Item(4, 20).Value2 = 0.000411479651105099

$ws.Cells.Item(5, 7).Value2 = 66.47695399999999
$ws.Cells.Item(5, 8).Value2 = 199.430862
$ws.Cells.Item(5, 9).Value2 = 0.04311983106164722
$ws.Cells.Item(5, 10).Value2 = 0.04311983106164721
$ws.Cells.Item(5, 11).Value2 = 3
$ws.Cells.Item(5, 12).Value2 = 1
$ws.Cells.Item(5, 13).Value2 = 4.073266333333334
$ws.Cells.Item(5, 14).Value2 = 12.219799
$ws.Cells.Item(5, 15).Value2 = 0.4378694427424867
$ws.Cells.Item(5, 16).Value2 = 0.4378694427424867
$ws.Cells.Item(5, 17).Value2 = 270.7783386707487
$ws.Cells.Item(5, 18).Value2 = 2437.005048036738
$ws.Cells.Item(5, 19).Value2 = 0.01888085639811364
$ws.Cells.Item(5, 20).Value2 = 0.01888085639811363

$ws.Cells.Item(6, 7).Value2 = 66.47695399999999
$ws.Cells.Item(6, 8).Value2 = 199.430862
$ws.Cells.Item(6, 9).Value2 = 0.04311983106164722
$ws.Cells.Item(6, 10).Value2 = 0.04311983106164721
$ws.Cells.Item(6, 11).Value2 = 3
$ws.Cells.Item(6, 12).Value2 = 1
$ws.Cells.Item(6, 13).Value2 = 0.1428486666666667
$ws.Cells.Item(6, 14).Value2 = 0.428546
$ws.Cells.Item(6, 15).Value2 = 0.01535599711660737
$ws.Cells.Item(6, 16).Value2 = 0.01535599711660737
$ws.Cells.Item(6, 17).Value2 = 9.496144242961332
$ws.Cells.Item(6, 18).Value2 = 85.46529818665199
$ws.Cells.Item(6, 19).Value2 = 0.0006621480014512518
$ws.Cells.Item(6, 20).Value2 = 0.0006621480014512517

$ws.Cells.Item(7, 7).Value2 = 1361.379069
$ws.Cells.Item(7, 8).Value2 = 4084.137207
$ws.Cells.Item(7, 9).Value2 = 0.8830494168872806
$ws.Cells.Item(7, 10).Value2 = 0.8830494168872804
$ws.Cells.Item(7, 11).Value2 = 3.0
$ws.Cells.Item(7, 12).Value2 = 1.0
$ws.Cells.Item(7, 13).Value2 = 0.5531836666666666
$ws.Cells.Item(7, 14).Value2 = 1.659551
$ws.Cells.Item(7, 15).Value2 = 0.05946633586794156
$ws.Cells.Item(7, 16).Value2 = 0.05946633586794157
$ws.Cells.Item(7, 17).Value2 = 753.0926651126729
$ws.Cells.Item(7, 18).Value2 = 6777.833986014057
$ws.Cells.Item(7, 19).Value2 = 0.05251171321260897
$ws.Cells.Item(7, 20).Value2 = 0.05251171321260897

$ws.Cells.Item(8, 7).Value2 = 1361.379069
$ws.Cells.Item(8, 8).Value2 = 4084.137207
$ws.Cells.Item(8, 9).Value2 = 0.8830494168872806
$ws.Cells.Item(8, 10).Value2 = 0.8830494168872804
$ws.Cells.Item(8, 11).Value2 = 3
$ws.Cells.Item(8, 12).Value2 = 1
$ws.Cells.Item(8, 13).Value2 = 4.444398333333333
$ws.Cells.Item(8, 14).Value2 = 13.333195
$ws.Cells.Item(8, 15).Value2 = 0.4777655233631019
$ws.Cells.Item(8, 16).Value2 = 0.4777655233631019
$ws.Cells.Item(8, 17).Value2 = 6050.510865298484
$ws.Cells.Item(8, 18).Value2 = 54454.59778768636
$ws.Cells.Item(8, 19).Value2 = 0.4218905668146335
$ws.Cells.Item(8, 20).Value2 = 0.4218905668146335

$ws.Cells.Item(9, 7).Value2 = 1361.379069
$ws.Cells.Item(9, 8).Value2 = 4084.137207
$ws.Cells.Item(9, 9).Value2 = 0.8830494168872806
$ws.Cells.Item(9, 10).Value2 = 0.8830494168872804
$ws.Cells.Item(9, 11).Value2 = 2
$ws.Cells.Item(9, 12).Value2 = 0.6666666666666666
$ws.Cells.Item(9, 13).Value2 = 0.08877066666666666
$ws.Cells.Item(9, 14).Value2 = 0.266312
$ws.Cells.Item(9, 15).Value2 = 0.009542700909862518
$ws.Cells.Item(9, 16).Value2 = 0.00954270090986252
$ws.Cells.Item(9, 17).Value2 = 120.850527541176
$ws.Cells.Item(9, 18).Value2 = 1087.654747870584
$ws.Cells.Item(9, 19).Value2 = 0.008426676473983818
$ws.Cells.Item(9, 20).Value2 = 0.008426676473983818

$ws.Cells.Item(10, 7).Value2 = 1361.379069
$ws.Cells.Item(10, 8).Value2 = 4084.137207
$ws.Cells.Item(10, 9).Value2 = 0.8830494168872806
$ws.Cells.Item(10, 10).Value2 = 0.8830494168872804
$ws.Cells.Item(10, 11).Value2 = 3
$ws.Cells.Item(10, 12).Value2 = 1
$ws.Cells.Item(10, 13).Value2 = 4.073266333333334
$ws.Cells.Item(10, 14).Value2 = 12.219799
$ws.Cells.Item(10, 15).Value2 = 0.4378694427424867
$ws.Cells.Item(10, 16).Value2 = 0.4378694427424867
$ws.Cells.Item(10, 17).Value2 = 5545.259528662376
$ws.Cells.Item(10, 18).Value2 = 49907.33575796139
$ws.Cells.Item(10, 19).Value2 = 0.3866603560865113
$ws.Cells.Item(10, 20).Value2 = 0.3866603560865113

$ws.Cells.Item(11, 7).Value2 = 1361.379069
$ws.Cells.Item(11, 8).Value2 = 4084.137207
$ws.Cells.Item(11, 9).Value2 = 0.8830494168872806
$ws.Cells.Item(11, 10).Value2 = 0.8830494168872804
$ws.Cells.Item(11, 11).Value2 = 3
$ws.Cells.Item(11, 12).Value2 = 1
$ws.Cells.Item(11, 13).Value2 = 0.1428486666666667
$ws.Cells.Item(11, 14).Value2 = 0.428546
$ws.Cells.Item(11, 15).Value2 = 0.01535599711660737
$ws.Cells.Item(11, 16).Value2 = 0.01535599711660737
$ws.Cells.Item(11, 17).Value2 = 194.471184834558
$ws.Cells.Item(11, 18).Value2 = 1750.240663511022
$ws.Cells.Item(11, 19).Value2 = 0.0135601042995429
$ws.Cells.Item(11, 20).Value2 = 0.0135601042995429

$ws.Cells.Item(12, 7).Value2 = 44.831112
$ws.Cells.Item(12, 8).Value2 = 134.493336
$ws.Cells.Item(12, 9).Value2 = 0.02907940059566787
$ws.Cells.Item(12, 10).Value2 = 0.02907940059566786
$ws.Cells.Item(12, 11).Value2 = 3.0
$ws.Cells.Item(12, 12).Value2 = 1.0
$ws.Cells.Item(12, 13).Value2 = 0.5531836666666666
$ws.Cells.Item(12, 14).Value2 = 1.659551
$ws.Cells.Item(12, 15).Value2 = 0.05946633586794156
$ws.Cells.Item(12, 16).Value2 = 0.05946633586794157
$ws.Cells.Item(12, 17).Value2 = 24.799838916904
$ws.Cells.Item(12, 18).Value2 = 223.198550252136
$ws.Cells.Item(12, 19).Value2 = 0.001729245402660406
$ws.Cells.Item(12, 20).Value2 = 0.001729245402660405

$ws.Cells.Item(13, 7).Value2 = 44.831112
$ws.Cells.Item(13, 8).Value2 = 134.493336
$ws.Cells.Item(13, 9).Value2 = 0.02907940059566787
$ws.Cells.Item(13, 10).Value2 = 0.02907940059566786
$ws.Cells.Item(13, 11).Value2 = 3
$ws.Cells.Item(13, 12).Value2 = 1
$ws.Cells.Item(13, 13).Value2 = 4.444398333333333
$ws.Cells.Item(13, 14).Value2 = 13.333195
$ws.Cells.Item(13, 15).Value2 = 0.4777655233631019
$ws.Cells.Item(13, 16).Value2 = 0.4777655233631019
$ws.Cells.Item(13, 17).Value2 = 199.24731945428
$ws.Cells.Item(13, 18).Value2 = 1793.22587508852
$ws.Cells.Item(13, 19).Value2 = 0.01389313504467456
$ws.Cells.Item(13, 20).Value2 = 0.01389313504467455

$ws.Cells.Item(14, 7).Value2 = 44.831112
$ws.Cells.Item(14, 8).Value2 = 134.493336
$ws.Cells.Item(14, 9).Value2 = 0.02907940059566787
$ws.Cells.Item(14, 10).Value2 = 0.02907940059566786
$ws.Cells.Item(14, 11).Value2 = 2
$ws.Cells.Item(14, 12).Value2 = 0.6666666666666666
$ws.Cells.Item(14, 13).Value2 = 0.08877066666666666
$ws.Cells.Item(14, 14).Value2 = 0.266312
$ws.Cells.Item(14, 15).Value2 = 0.009542700909862518
$ws.Cells.Item(14, 16).Value2 = 0.00954270090986252
$ws.Cells.Item(14, 17).Value2 = 3.979687699648
$ws.Cells.Item(14, 18).Value2 = 35.817189296832
$ws.Cells.Item(14, 19).Value2 = 0.0002774960225225364
$ws.Cells.Item(14, 20).Value2 = 0.0002774960225225364

$ws.Cells.Item(15, 7).Value2 = 44.831112
$ws.Cells.Item(15, 8).Value2 = 134.493336
$ws.Cells.Item(15, 9).Value2 = 0.02907940059566787
$ws.Cells.Item(15, 10).Value2 = 0.02907940059566786
$ws.Cells.Item(15, 11).Value2 = 3
$ws.Cells.Item(15, 12).Value2 = 1
$ws.Cells.Item(15, 13).Value2 = 4.073266333333334
$ws.Cells.Item(15, 14).Value2 = 12.219799
$ws.Cells.Item(15, 15).Value2 = 0.4378694427424867
$ws.Cells.Item(15, 16).Value2 = 0.4378694427424867
$ws.Cells.Item(15, 17).Value2 = 182.609059195496
$ws.Cells.Item(15, 18).Value2 = 1643.481532759464
$ws.Cells.Item(15, 19).Value2 = 0.01273298093411063
$ws.Cells.Item(15, 20).Value2 = 0.01273298093411062

$ws.Cells.Item(16, 7).Value2 = 44.831112
$ws.Cells.Item(16, 8).Value2 = 134.493336
$ws.Cells.Item(16, 9).Value2 = 0.02907940059566787
$ws.Cells.Item(16, 10).Value2 = 0.02907940059566786
$ws.Cells.Item(16, 11).Value2 = 3
$ws.Cells.Item(16, 12).Value2 = 1
$ws.Cells.Item(16, 13).Value2 = 0.1428486666666667
$ws.Cells.Item(16, 14).Value2 = 0.428546
$ws.Cells.Item(16, 15).Value2 = 0.01535599711660737
$ws.Cells.Item(16, 16).Value2 = 0.01535599711660737
$ws.Cells.Item(16, 17).Value2 = 6.404064574383999
$ws.Cells.Item(16, 18).Value2 = 57.636581169456
$ws.Cells.Item(16, 19).Value2 = 0.0004465431916997466
$ws.Cells.Item(16, 20).Value2 = 0.0004465431916997465

$ws.Cells.Item(17, 7).Value2 = 52.83062100000001
$ws.Cells.Item(17, 8).Value2 = 158.491863
$ws.Cells.Item(17, 9).Value2 = 0.0342682285413064
$ws.Cells.Item(17, 10).Value2 = 0.03426822854130639
$ws.Cells.Item(17, 11).Value2 = 3.0
$ws.Cells.Item(17, 12).Value2 = 1.0
$ws.Cells.Item(17, 13).Value2 = 0.5531836666666666
$ws.Cells.Item(17, 14).Value2 = 1.659551
$ws.Cells.Item(17, 15).Value2 = 0.05946633586794156
$ws.Cells.Item(17, 16).Value2 = 0.05946633586794157
$ws.Cells.Item(17, 17).Value2 = 29.225036637057
$ws.Cells.Item(17, 18).Value2 = 263.025329733513
$ws.Cells.Item(17, 19).Value2 = 0.002037805988036707
$ws.Cells.Item(17, 20).Value2 = 0.002037805988036707

$ws.Cells.Item(18, 7).Value2 = 52.83062100000001
$ws.Cells.Item(18, 8).Value2 = 158.491863
$ws.Cells.Item(18, 9).Value2 = 0.0342682285413064
$ws.Cells.Item(18, 10).Value2 = 0.03426822854130639
$ws.Cells.Item(18, 11).Value2 = 3
$ws.Cells.Item(18, 12).Value2 = 1
$ws.Cells.Item(18, 13).Value2 = 4.444398333333333
$ws.Cells.Item(18, 14).Value2 = 13.333195
$ws.Cells.Item(18, 15).Value2 = 0.4777655233631019
$ws.Cells.Item(18, 16).Value2 = 0.4777655233631019
$ws.Cells.Item(18, 17).Value2 = 234.800323921365
$ws.Cells.Item(18, 18).Value2 = 2113.202915292286
$ws.Cells.Item(18, 19).Value2 = 0.01637217814376364
$ws.Cells.Item(18, 20).Value2 = 0.01637217814376363

$ws.Cells.Item(19, 7).Value2 = 52.83062100000001
$ws.Cells.Item(19, 8).Value2 = 158.491863
$ws.Cells.Item(19, 9).Value2 = 0.0342682285413064
$ws.Cells.Item(19, 10).Value2 = 0.03426822854130639
$ws.Cells.Item(19, 11).Value2 = 2
$ws.Cells.Item(19, 12).Value2 = 0.6666666666666666
$ws.Cells.Item(19, 13).Value2 = 0.08877066666666666
$ws.Cells.Item(19, 14).Value2 = 0.266312
$ws.Cells.Item(19, 15).Value2 = 0.009542700909862518
$ws.Cells.Item(19, 16).Value2 = 0.00954270090986252
$ws.Cells.Item(19, 17).Value2 = 4.689809446584
$ws.Cells.Item(19, 18).Value2 = 42.208285019256
$ws.Cells.Item(19, 19).Value2 = 0.0003270114556805013
$ws.Cells.Item(19, 20).Value2 = 0.0003270114556805013

$ws.Cells.Item(20, 7).Value2 = 52.83062100000001
$ws.Cells.Item(20, 8).Value2 = 158.491863
$ws.Cells.Item(20, 9).Value2 = 0.0342682285413064
$ws.Cells.Item(20, 10).Value2 = 0.03426822854130639
$ws.Cells.Item(20, 11).Value2 = 3
$ws.Cells.Item(20, 12).Value2 = 1
$ws.Cells.Item(20, 13).Value2 = 4.073266333333334
$ws.Cells.Item(20, 14).Value2 = 12.219799
$ws.Cells.Item(20, 15).Value2 = 0.4378694427424867
$ws.Cells.Item(20, 16).Value2 = 0.4378694427424867
$ws.Cells.Item(20, 17).Value2 = 215.1931898883931
$ws.Cells.Item(20, 18).Value2 = 1936.738708995537
$ws.Cells.Item(20, 19).Value2 = 0.01500501013515401
$ws.Cells.Item(20, 20).Value2 = 0.01500501013515401

$ws.Cells.Item(21, 7).Value2 = 52.83062100000001
$ws.Cells.Item(21, 8).Value2 = 158.491863
$ws.Cells.Item(21, 9).Value2 = 0.0342682285413064
$ws.Cells.Item(21, 10).Value2 = 0.03426822854130639
$ws.Cells.Item(21, 11).Value2 = 3
$ws.Cells.Item(21, 12).Value2 = 1
$ws.Cells.Item(21, 13).Value2 = 0.1428486666666667
$ws.Cells.Item(21, 14).Value2 = 0.428546
$ws.Cells.Item(21, 15).Value2 = 0.01535599711660737
$ws.Cells.Item(21, 16).Value2 = 0.01535599711660737
$ws.Cells.Item(21, 17).Value2 = 7.546783769022
$ws.Cells.Item(21, 18).Value2 = 67.92105392119801
$ws.Cells.Item(21, 19).Value2 = 0.0005262228186715436
$ws.Cells.Item(21, 20).Value2 = 0.0005262228186715436

$ws.Cells.Item(22, 7).Value2 = 16.16161433333333
$ws.Cells.Item(22, 8).Value2 = 48.484843
$ws.Cells.Item(22, 9).Value2 = 0.01048312291409786
$ws.Cells.Item(22, 10).Value2 = 0.01048312291409786
$ws.Cells.Item(22, 11).Value2 = 3.0
$ws.Cells.Item(22, 12).Value2 = 1.0
$ws.Cells.Item(22, 13).Value2 = 0.5531836666666666
$ws.Cells.Item(22, 14).Value2 = 1.659551
$ws.Cells.Item(22, 15).Value2 = 0.05946633586794156
$ws.Cells.Item(22, 16).Value2 = 0.05946633586794157
$ws.Cells.Item(22, 17).Value2 = 8.940341076165888
$ws.Cells.Item(22, 18).Value2 = 80.463069685493
$ws.Cells.Item(22, 19).Value2 = 0.0006233929081546579
$ws.Cells.Item(22, 20).Value2 = 0.0006233929081546579

$ws.Cells.Item(23, 7).Value2 = 16.16161433333333
$ws.Cells.Item(23, 8).Value2 = 48.484843
$ws.Cells.Item(23, 9).Value2 = 0.01048312291409786
$ws.Cells.Item(23, 10).Value2 = 0.01048312291409786
$ws.Cells.Item(23, 11).Value2 = 3
$ws.Cells.Item(23, 12).Value2 = 1
$ws.Cells.Item(23, 13).Value2 = 4.444398333333333
$ws.Cells.Item(23, 14).Value2 = 13.333195
$ws.Cells.Item(23, 15).Value2 = 0.4777655233631019
$ws.Cells.Item(23, 16).Value2 = 0.4777655233631019
$ws.Cells.Item(23, 17).Value2 = 71.82865180704277
$ws.Cells.Item(23, 18).Value2 = 646.457866263385
$ws.Cells.Item(23, 19).Value2 = 0.005008474705533692
$ws.Cells.Item(23, 20).Value2 = 0.005008474705533691

$ws.Cells.Item(24, 7).Value2 = 16.16161433333333
$ws.Cells.Item(24, 8).Value2 = 48.484843
$ws.Cells.Item(24, 9).Value2 = 0.01048312291409786
$ws.Cells.Item(24, 10).Value2 = 0.01048312291409786
$ws.Cells.Item(24, 11).Value2 = 2
$ws.Cells.Item(24, 12).Value2 = 0.6666666666666666
$ws.Cells.Item(24, 13).Value2 = 0.08877066666666666
$ws.Cells.Item(24, 14).Value2 = 0.266312
$ws.Cells.Item(24, 15).Value2 = 0.009542700909862518
$ws.Cells.Item(24, 16).Value2 = 0.00954270090986252
$ws.Cells.Item(24, 17).Value2 = 1.434677278779555
$ws.Cells.Item(24, 18).Value2 = 12.912095509016
$ws.Cells.Item(24, 19).Value2 = 0.0001000373065705623
$ws.Cells.Item(24, 20).Value2 = 0.0001000373065705623

$ws.Cells.Item(25, 7).Value2 = 16.16161433333333
$ws.Cells.Item(25, 8).Value2 = 48.484843
$ws.Cells.Item(25, 9).Value2 = 0.01048312291409786
$ws.Cells.Item(25, 10).Value2 = 0.01048312291409786
$ws.Cells.Item(25, 11).Value2 = 3
$ws.Cells.Item(25, 12).Value2 = 1
$ws.Cells.Item(25, 13).Value2 = 4.073266333333334
$ws.Cells.Item(25, 14).Value2 = 12.219799
$ws.Cells.Item(25, 15).Value2 = 0.4378694427424867
$ws.Cells.Item(25, 16).Value2 = 0.4378694427424867
$ws.Cells.Item(25, 17).Value2 = 65.83055955628411
$ws.Cells.Item(25, 18).Value2 = 592.4750360065569
$ws.Cells.Item(25, 19).Value2 = 0.004590239188597025
$ws.Cells.Item(25, 20).Value2 = 0.004590239188597025

$ws.Cells.Item(26, 7).Value2 = 16.16161433333333
$ws.Cells.Item(26, 8).Value2 = 48.484843
$ws.Cells.Item(26, 9).Value2 = 0.01048312291409786
$ws.Cells.Item(26, 10).Value2 = 0.01048312291409786
$ws.Cells.Item(26, 11).Value2 = 3
$ws.Cells.Item(26, 12).Value2 = 1
$ws.Cells.Item(26, 13).Value2 = 0.1428486666666667
$ws.Cells.Item(26, 14).Value2 = 0.428546
$ws.Cells.Item(26, 15).Value2 = 0.01535599711660737
$ws.Cells.Item(26, 16).Value2 = 0.01535599711660737
$ws.Cells.Item(26, 17).Value2 = 2.308665058697555
$ws.Cells.Item(26, 18).Value2 = 20.777985528278
$ws.Cells.Item(26, 19).Value2 = 0.0001609788052419275
$ws.Cells.Item(26, 20).Value2 = 0.0001609788052419275

